$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quarter header labels (rows 8, 17, 27, 37) - shift forward by one quarter
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل چهارم منتهی به 1401/12"

$ws.Range("E17").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F17").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G17").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H17").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I17").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J17").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K17").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L17").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M17").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N17").Value = "فصل چهارم منتهی به 1401/12"

$ws.Range("E27").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F27").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G27").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H27").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I27").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J27").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K27").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L27").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M27").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N27").Value = "فصل چهارم منتهی به 1401/12"

$ws.Range("E37").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F37").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G37").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H37").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I37").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J37").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K37").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L37").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M37").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N37").Value = "فصل چهارم منتهی به 1401/12"

# Update data rows with refreshed quarterly figures
# Row 10
$ws.Range("E10").Value = 30
$ws.Range("F10").Value = 75
$ws.Range("G10").Value = 67
$ws.Range("H10").Value = 103
$ws.Range("I10").Value = 47
$ws.Range("J10").Value = "-"
$ws.Range("K10").Value = 48
$ws.Range("L10").Value = 53
$ws.Range("M10").Value = 37
$ws.Range("N10").Value = 21

# Row 11
$ws.Range("E11").Value = 8252
$ws.Range("F11").Value = 8588
$ws.Range("G11").Value = 9018
$ws.Range("H11").Value = 8913
$ws.Range("I11").Value = 8231
$ws.Range("J11").Value = "-"
$ws.Range("K11").Value = 15931
$ws.Range("L11").Value = 9343
$ws.Range("M11").Value = 11852
$ws.Range("N11").Value = 11740

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "-"
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0

# Row 13
$ws.Range("E13").Value = 8284
$ws.Range("F13").Value = 8664
$ws.Range("G13").Value = 9086
$ws.Range("H13").Value = 9016
$ws.Range("I13").Value = 8278
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 15979
$ws.Range("L13").Value = 9396
$ws.Range("M13").Value = 11889
$ws.Range("N13").Value = 11761

# Row 19
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = 87
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = "-"
$ws.Range("K19").Value = 51
$ws.Range("L19").Value = 56
$ws.Range("M19").Value = 36
$ws.Range("N19").Value = 26

# Row 20
$ws.Range("E20").Value = 10217
$ws.Range("F20").Value = 8352
$ws.Range("G20").Value = 10428
$ws.Range("H20").Value = 9300
$ws.Range("I20").Value = 7741
$ws.Range("J20").Value = "-"
$ws.Range("K20").Value = 14024
$ws.Range("L20").Value = 8517
$ws.Range("M20").Value = 11059
$ws.Range("N20").Value = 10427

# Row 21
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = "-"
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0

# Row 23
$ws.Range("E23").Value = 10235
$ws.Range("F23").Value = 8429
$ws.Range("G23").Value = 10497
$ws.Range("H23").Value = 9387
$ws.Range("I23").Value = 7784
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 14075
$ws.Range("L23").Value = 8573
$ws.Range("M23").Value = 11095
$ws.Range("N23").Value = 10453

# Row 29
$ws.Range("E29").Value = 17060
$ws.Range("F29").Value = 83838
$ws.Range("G29").Value = 81240
$ws.Range("H29").Value = 112068
$ws.Range("I29").Value = 59818
$ws.Range("J29").Value = "-"
$ws.Range("K29").Value = 72477
$ws.Range("L29").Value = 85231
$ws.Range("M29").Value = 58136
$ws.Range("N29").Value = 50989

# Row 30
$ws.Range("E30").Value = 568266
$ws.Range("F30").Value = 635275
$ws.Range("G30").Value = 785259
$ws.Range("H30").Value = 769516
$ws.Range("I30").Value = 785184
$ws.Range("J30").Value = "-"
$ws.Range("K30").Value = 1277764
$ws.Range("L30").Value = 981666
$ws.Range("M30").Value = 1015218
$ws.Range("N30").Value = 1492067

# Row 31
$ws.Range("E31").Value = 5334
$ws.Range("F31").Value = 2057
$ws.Range("G31").Value = 295
$ws.Range("H31").Value = 2775
$ws.Range("I31").Value = 269
$ws.Range("J31").Value = "-"
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0

# Row 33
$ws.Range("E33").Value = 590660
$ws.Range("F33").Value = 721170
$ws.Range("G33").Value = 866794
$ws.Range("H33").Value = 884359
$ws.Range("I33").Value = 845271
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1350241
$ws.Range("L33").Value = 1066897
$ws.Range("M33").Value = 1073354
$ws.Range("N33").Value = 1543056

# Row 39
$ws.Range("E39").Value = 1066250000
$ws.Range("F39").Value = 1117840000
$ws.Range("G39").Value = 1194705882
$ws.Range("H39").Value = 1288137931
$ws.Range("I39").Value = 1391372093
$ws.Range("J39").Value = 1309000000
$ws.Range("K39").Value = 1421117647
$ws.Range("L39").Value = 1521982143
$ws.Range("M39").Value = 1611342857
$ws.Range("N39").Value = 1961115385

# Row 40
$ws.Range("E40").Value = 55221595
$ws.Range("F40").Value = 76062620
$ws.Range("G40").Value = 75302934
$ws.Range("H40").Value = 82743656
$ws.Range("I40").Value = 102613884
$ws.Range("J40").Value = 89563907
$ws.Range("K40").Value = 91112664
$ws.Range("L40").Value = 115259598
$ws.Range("M40").Value = 96254573
$ws.Range("N40").Value = 143096480

# Row 41
$ws.Range("E41").Value = 384500000
$ws.Range("F41").Value = 1028500000
$ws.Range("G41").Value = 295000000
$ws.Range("H41").Value = "-"
$ws.Range("I41").Value = "-"
$ws.Range("J41").Value = "-"
$ws.Range("K41").Value = "-"
$ws.Range("L41").Value = "-"
$ws.Range("M41").Value = "-"
$ws.Range("N41").Value = "-"

